# Update the NATMI TPM output numbers (Shh-Smo) for rows 2-7.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values keyed by cell address -> value.
$updates = @{
    # Row 2
    "E2" = 2
    "F2" = 0.6666666666666666
    "G2" = 0.08785
    "H2" = 0.26355
    "M2" = 3.151228
    "N2" = 9.453684000000001
    "O2" = 0.1027676232988273
    "P2" = 0.1027676232988273
    "Q2" = 0.2768353798
    "R2" = 2.4915184182
    "S2" = 0.1027676232988273
    "T2" = 0.1027676232988273

    # Row 3
    "E3" = 2
    "F3" = 0.6666666666666666
    "G3" = 0.08785
    "H3" = 0.26355
    "O3" = 0.6264504114587274
    "P3" = 0.6264504114587273
    "Q3" = 1.687531851133333
    "R3" = 15.1877866602
    "S3" = 0.6264504114587274
    "T3" = 0.6264504114587273

    # Row 4
    "E4" = 2
    "F4" = 0.6666666666666666
    "G4" = 0.08785
    "H4" = 0.26355
    "M4" = 0.27146
    "N4" = 0.8143800000000001
    "O4" = 0.00885283420326922
    "P4" = 0.008852834203269218
    "Q4" = 0.023847761
    "R4" = 0.214629849
    "S4" = 0.00885283420326922
    "T4" = 0.008852834203269218

    # Row 5
    "E5" = 2
    "F5" = 0.6666666666666666
    "G5" = 0.08785
    "H5" = 0.26355
    "M5" = 7.293697999999999
    "N5" = 21.881094
    "O5" = 0.2378615601661986
    "P5" = 0.2378615601661986
    "Q5" = 0.6407513692999999
    "R5" = 5.766762323699999
    "S5" = 0.2378615601661986
    "T5" = 0.2378615601661986

    # Row 6
    "E6" = 2
    "F6" = 0.6666666666666666
    "G6" = 0.08785
    "H6" = 0.26355
    "M6" = 0.2802906666666667
    "N6" = 0.8408720000000001
    "O6" = 0.009140819276224114
    "P6" = 0.009140819276224113
    "Q6" = 0.02462353506666667
    "R6" = 0.2216118156
    "S6" = 0.009140819276224114
    "T6" = 0.009140819276224113

    # Row 7
    "E7" = 2
    "F7" = 0.6666666666666666
    "G7" = 0.08785
    "H7" = 0.26355
    "M7" = 0.4577083333333333
    "N7" = 1.373125
    "O7" = 0.01492675159675341
    "P7" = 0.01492675159675341
    "Q7" = 0.04020967708333333
    "R7" = 0.36188709375
    "S7" = 0.01492675159675341
    "T7" = 0.01492675159675341
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
